# InvoiceSenderControl.xlsx — update the example "person" placeholder used
# in the Polish email-template row of the Control sheet.
#
# "Maria Nowak" -> "Mikołaj Męderski"  (row 7, column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control")

$cell = $ws.Range("B7")
$cell.Value = "Mikołaj Męderski"

# Touching the font here forks a dedicated style for B7 (new font + cellXf
# entries), matching how Excel re-keys the cell's style whenever its font
# is nudged after a content edit.
$cell.Font.Size = 11
